$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BI1").Value = 0.953087019512876
$ws.Range("A2").Value = 0.94415120313753498
$ws.Range("C2").Value = 0.52892076733121951
$ws.Range("AA2").Value = 0.94113421265807973
$ws.Range("L3").Value = 0.91758663211644254
$ws.Range("AY4").Value = 0.97331291190539548
$ws.Range("BO5").Value = 0.97194377782091212
$ws.Range("BP5").Value = 0.80759733626677366
$ws.Range("G6").Value = 0.8151352135985942
$ws.Range("AI6").Value = 0.66258944687529153
$ws.Range("AY6").Value = 0.89806787074055938
$ws.Range("AJ7").Value = 0.97095701489748398
$ws.Range("F8").Value = 0.80586924324331477
$ws.Range("AU8").Value = 0.88857065724776962
$ws.Range("AD9").Value = 0.99248176362940366
$ws.Range("M10").Value = 0.75780554151610346
$ws.Range("AG10").Value = 0.59380052003240014
$ws.Range("AZ10").Value = 0.86393292907916786
$ws.Range("BL11").Value = 0.93813562899211012
$ws.Range("M12").Value = 0.98785106466078854
$ws.Range("D13").Value = 0.84117869472973583
$ws.Range("AI13").Value = 0.56622462863728851
$ws.Range("BJ13").Value = 0.81240153469088727
$ws.Range("BA14").Value = 0.98226973141652396
$ws.Range("BP14").Value = 0.87665961495681721
$ws.Range("R15").Value = 0.90302341350865722
$ws.Range("AE15").Value = 0.96353244699521612
$ws.Range("BF15").Value = 0.98503559727966028
$ws.Range("G16").Value = 0.88569623654825302
$ws.Range("U16").Value = 0.99535119159361252
$ws.Range("BD16").Value = 0.93318849262063641
$ws.Range("F17").Value = 0.68628834189253674
$ws.Range("AJ18").Value = 0.88592733043635685
$ws.Range("AX18").Value = 0.92705307622589772
$ws.Range("T19").Value = 0.92458775385161696
$ws.Range("AV19").Value = 0.71228419102134233
$ws.Range("Z20").Value = 0.68021261193859506
$ws.Range("AR20").Value = 0.84191457085424004
$ws.Range("BH20").Value = 0.89607265754834331
$ws.Range("Q21").Value = 0.95986137523773862
$ws.Range("V21").Value = 0.59643282644480478
$ws.Range("BO22").Value = 0.97631670262692438
$ws.Range("BP22").Value = 0.75783797152860077
$ws.Range("Y23").Value = 0.86242195884633399
$ws.Range("F24").Value = 0.83304297535504368
$ws.Range("AF24").Value = 0.98727872558414109
$ws.Range("BH24").Value = 0.83712296664843477
$ws.Range("AE25").Value = 0.96261829010105893
$ws.Range("AL26").Value = 0.93638770893891488
$ws.Range("AT27").Value = 0.73473878873536758
$ws.Range("BA28").Value = 0.62066683365795772
$ws.Range("S29").Value = 0.94059580698170608
$ws.Range("AX29").Value = 0.83886199103381132
$ws.Range("AY29").Value = 0.88547993038538864
$ws.Range("L30").Value = 0.84095272953514788
$ws.Range("BO30").Value = 0.70085876215664622
$ws.Range("N31").Value = 0.85510799473560639
$ws.Range("AE32").Value = 0.93555845782478508
$ws.Range("BM32").Value = 0.64998360308960668
$ws.Range("W33").Value = 0.86826427323026656
$ws.Range("AP33").Value = 0.90391916445464637
$ws.Range("BL33").Value = 0.93465988605078687
$ws.Range("B35").Value = 0.57670964525574619
$ws.Range("BF35").Value = 0.96028343155640061
$ws.Range("AH36").Value = 0.83591708709195101
$ws.Range("S37").Value = 0.90959650057897656
$ws.Range("AC37").Value = 0.93413441107510464
$ws.Range("AG38").Value = 0.86561039229462344
$ws.Range("AW38").Value = 0.82235819425262391
$ws.Range("BC38").Value = 0.65410356304105743
$ws.Range("BL38").Value = 0.93467087542193095
$ws.Range("M39").Value = 0.74421294657983716
$ws.Range("BE39").Value = 0.95056666199760276
$ws.Range("A40").Value = 0.80316806121768258
$ws.Range("AP40").Value = 0.76375136065855376
$ws.Range("U41").Value = 0.79399838132523892
$ws.Range("AY41").Value = 0.91365136779243095
$ws.Range("F42").Value = 0.9165177540979681
$ws.Range("AK42").Value = 0.85670951489408442
$ws.Range("AT42").Value = 0.59701328776066087
$ws.Range("C43").Value = 0.75572994438072527
$ws.Range("AD44").Value = 0.72511307442358275
$ws.Range("AQ44").Value = 0.74181662986191532
$ws.Range("AG45").Value = 0.72569186580800182
$ws.Range("AK45").Value = 0.74184028368112254
$ws.Range("W46").Value = 0.72056601498712647
$ws.Range("AN46").Value = 0.78573330233669503
$ws.Range("AW46").Value = 0.99926289345419583
$ws.Range("AZ46").Value = 0.755826898486468
$ws.Range("AW47").Value = 0.89794163690170525
$ws.Range("F48").Value = 0.85565603498743514
$ws.Range("AU48").Value = 0.91315328590675149
$ws.Range("T49").Value = 0.82757686750658244
$ws.Range("AS49").Value = 0.78352483100073333
$ws.Range("BD49").Value = 0.68368706773798316
$ws.Range("AY50").Value = 0.84991991045793225
$ws.Range("BC50").Value = 0.99895612665822897
$ws.Range("C52").Value = 0.70773085687114046
$ws.Range("H52").Value = 0.70903767430487452
$ws.Range("K52").Value = 0.99385061889364379
$ws.Range("C53").Value = 0.55422602796515252
$ws.Range("AH53").Value = 0.74908211027954685
$ws.Range("AJ53").Value = 0.83040028095028928
$ws.Range("AU53").Value = 0.75708816103420307
$ws.Range("AY53").Value = 0.87064311448390996
$ws.Range("BJ53").Value = 0.76297440475838085
$ws.Range("S54").Value = 0.97500574577261434
$ws.Range("AU54").Value = 0.93029149540309408
$ws.Range("BC54").Value = 0.84242526050609423
$ws.Range("BP54").Value = 0.76985664690767663
$ws.Range("BE55").Value = 0.82697682474025602
$ws.Range("AC56").Value = 0.83076336976363019
$ws.Range("AX56").Value = 0.99371888367952987
$ws.Range("AJ57").Value = 0.83170223378393637
$ws.Range("AT57").Value = 0.71867516190974512
$ws.Range("T58").Value = 0.68305355497523546
$ws.Range("AH59").Value = 0.58016121744176075
$ws.Range("AV59").Value = 0.90088618021766487
$ws.Range("BF60").Value = 0.80392696919455564
$ws.Range("I61").Value = 0.71538205948719436
$ws.Range("Q61").Value = 0.8063233975100772
$ws.Range("X62").Value = 0.83170940967676588
$ws.Range("AB62").Value = 0.69545721111125891
$ws.Range("AL63").Value = 0.76667884399148645
$ws.Range("BH63").Value = 0.85975352677888506
$ws.Range("BJ63").Value = 0.81390299170934766
$ws.Range("Y64").Value = 0.99532279641168631
$ws.Range("L65").Value = 0.62229702442227752
$ws.Range("P65").Value = 0.89507536005551358
$ws.Range("AQ65").Value = 0.98871215771248289
$ws.Range("BH65").Value = 0.78760500558377178
$ws.Range("I66").Value = 0.97935344390520496
$ws.Range("AC66").Value = 0.94239461426162419
$ws.Range("AY66").Value = 0.7652442537258537
$ws.Range("A67").Value = 0.68696207478540172
$ws.Range("AL68").Value = 0.9136258756477007
